# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets, reflecting the data refresh from
# the generator run at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 333
    $ws.Range("F3").Value = 1349
    $ws.Range("F5").Value = 70
}
